$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H12").Value = 877.5
$ws.Range("I12").Value = 855
$ws.Range("K12").Value = 855
$ws.Range("M12").Value = -685
$ws.Range("H76").Value = 8778.799999999999
$ws.Range("I76").Value = 6697.25
$ws.Range("K76").Value = 6697.25
$ws.Range("M76").Value = -6382.25
$ws.Range("H79").Value = 8778.799999999999
$ws.Range("I79").Value = 6697.25
$ws.Range("K79").Value = 6697.25
$ws.Range("M79").Value = -5605.25
$ws.Range("H94").Value = 628.5
$ws.Range("I94").Value = 338
$ws.Range("K94").Value = 338
$ws.Range("M94").Value = 113
$ws.Range("H98").Value = 2319
$ws.Range("I98").Value = 1204.4166
$ws.Range("J98").Value = 5662.75
$ws.Range("K98").Value = 1204.4166
$ws.Range("L98").Value = 5662.75
$ws.Range("M98").Value = 293.5834
$ws.Range("N98").Value = -8658.75
$ws.Range("H101").Value = 382.46155
$ws.Range("I101").Value = 412.54544
$ws.Range("J101").Value = 217
$ws.Range("K101").Value = 1237.63632
$ws.Range("L101").Value = 651
$ws.Range("M101").Value = 384.3636799999999
$ws.Range("N101").Value = -3895
$ws.Range("H112").Value = 2339.5715
$ws.Range("I112").Value = 1597.8
$ws.Range("K112").Value = 4793.4
$ws.Range("M112").Value = -3685.4
$ws.Range("H122").Value = 2319
$ws.Range("I122").Value = 1204.4166
$ws.Range("J122").Value = 5662.75
$ws.Range("K122").Value = 3613.2498
$ws.Range("L122").Value = 16988.25
$ws.Range("M122").Value = -1163.2498
$ws.Range("N122").Value = -21888.25
$ws.Range("H138").Value = 3487.425
$ws.Range("I138").Value = 2822.8572
$ws.Range("J138").Value = 4221.9473
$ws.Range("K138").Value = 8468.571599999999
$ws.Range("L138").Value = 12665.8419
$ws.Range("M138").Value = -3328.571599999999
$ws.Range("N138").Value = -22945.8419

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 3814.1833
$ws.Range("I32").Value = 3402.66
$ws.Range("J32").Value = 5871.8
$ws.Range("K32").Value = 3402.66
$ws.Range("L32").Value = 5871.8
$ws.Range("M32").Value = -3115.66
$ws.Range("N32").Value = -6445.8
$ws.Range("H122").Value = 2222.743
$ws.Range("I122").Value = 1747.2667
$ws.Range("K122").Value = 5241.800099999999
$ws.Range("M122").Value = -2791.800099999999
$ws.Range("H132").Value = 2800.0715
$ws.Range("I132").Value = 1498.375
$ws.Range("K132").Value = 4495.125
$ws.Range("M132").Value = -1965.125

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 2594.6875
$ws.Range("I134").Value = 1572.7142
$ws.Range("K134").Value = 4718.142599999999
$ws.Range("M134").Value = -2183.142599999999

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 1700.1666
$ws.Range("I22").Value = 625.9091
$ws.Range("J22").Value = 3388.2856
$ws.Range("K22").Value = 625.9091
$ws.Range("L22").Value = 3388.2856
$ws.Range("M22").Value = -275.9091
$ws.Range("N22").Value = -4088.2856
$ws.Range("H60").Value = 81984
$ws.Range("I60").Value = 50000
$ws.Range("J60").Value = 89980
$ws.Range("K60").Value = 50000
$ws.Range("L60").Value = 89980
$ws.Range("M60").Value = -49489
$ws.Range("N60").Value = -91002
$ws.Range("H74").Value = 134031.25
$ws.Range("I74").Value = 24166.666
$ws.Range("J74").Value = 199950
$ws.Range("K74").Value = 24166.666
$ws.Range("L74").Value = 199950
$ws.Range("M74").Value = -23292.666
$ws.Range("N74").Value = -201698
$ws.Range("H77").Value = 134031.25
$ws.Range("I77").Value = 24166.666
$ws.Range("J77").Value = 199950
$ws.Range("K77").Value = 72499.99800000001
$ws.Range("L77").Value = 599850
$ws.Range("M77").Value = -68131.99800000001
$ws.Range("N77").Value = -608586
$ws.Range("H105").Value = 2955.5356
$ws.Range("I105").Value = 2605
$ws.Range("K105").Value = 2605
$ws.Range("M105").Value = -858

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H55").Value = 3251.875
$ws.Range("J55").Value = 3251.875
$ws.Range("L55").Value = 9755.625
$ws.Range("N55").Value = -10109.625
$ws.Range("H131").Value = 1400734.8
$ws.Range("I131").Value = 933
$ws.Range("J131").Value = 1820675.2
$ws.Range("K131").Value = 2799
$ws.Range("L131").Value = 5462025.6
$ws.Range("M131").Value = 2241
$ws.Range("N131").Value = -5472105.6
$ws.Range("H136").Value = 2563.8
$ws.Range("I136").Value = 1706.4445
$ws.Range("J136").Value = 3849.8333
$ws.Range("K136").Value = 5119.333500000001
$ws.Range("L136").Value = 11549.4999
$ws.Range("M136").Value = -19.33350000000064
$ws.Range("N136").Value = -21749.4999

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H46").Value = 47499.4
$ws.Range("J46").Value = 47499.4
$ws.Range("L46").Value = 47499.4
$ws.Range("N46").Value = -47811.4
$ws.Range("H97").Value = 472.96
$ws.Range("I97").Value = 465.26315
$ws.Range("J97").Value = 497.33334
$ws.Range("K97").Value = 465.26315
$ws.Range("L97").Value = 497.33334
$ws.Range("M97").Value = 30.73685
$ws.Range("N97").Value = -1489.33334
$ws.Range("H102").Value = 2242.8386
$ws.Range("I102").Value = 1447.7407
$ws.Range("J102").Value = 7609.75
$ws.Range("K102").Value = 1447.7407
$ws.Range("L102").Value = 7609.75
$ws.Range("M102").Value = 174.2592999999999
$ws.Range("N102").Value = -10853.75
$ws.Range("H122").Value = 3031.3547
$ws.Range("I122").Value = 1424.3636
$ws.Range("J122").Value = 6959.5557
$ws.Range("K122").Value = 4273.0908
$ws.Range("L122").Value = 20878.6671
$ws.Range("M122").Value = -1823.0908
$ws.Range("N122").Value = -25778.6671
$ws.Range("H132").Value = 3418.3
$ws.Range("I132").Value = 2796.875
$ws.Range("K132").Value = 8390.625
$ws.Range("M132").Value = -5860.625

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 2035.1305
$ws.Range("I16").Value = 1291.8
$ws.Range("J16").Value = 6990.6665
$ws.Range("K16").Value = 1291.8
$ws.Range("L16").Value = 6990.6665
$ws.Range("M16").Value = -1121.8
$ws.Range("N16").Value = -7330.6665
$ws.Range("H46").Value = 3816.9412
$ws.Range("I46").Value = 2027.1428
$ws.Range("J46").Value = 5069.8
$ws.Range("K46").Value = 2027.1428
$ws.Range("L46").Value = 5069.8
$ws.Range("M46").Value = -1839.1428
$ws.Range("N46").Value = -5445.8
$ws.Range("H122").Value = 6474.75
$ws.Range("I122").Value = 6474.75
$ws.Range("K122").Value = 19424.25
$ws.Range("M122").Value = -16974.25
$ws.Range("H132").Value = 12548.565
$ws.Range("I132").Value = 8971.154
$ws.Range("J132").Value = 17199.2
$ws.Range("K132").Value = 26913.462
$ws.Range("L132").Value = 51597.60000000001
$ws.Range("M132").Value = -24383.462
$ws.Range("N132").Value = -56657.60000000001

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H42").Value = 68450
$ws.Range("J42").Value = 87900
$ws.Range("L42").Value = 87900
$ws.Range("N42").Value = -88656
$ws.Range("H96").Value = 2046.1111
$ws.Range("I96").Value = 694.5
$ws.Range("K96").Value = 694.5
$ws.Range("M96").Value = 678.5
$ws.Range("H136").Value = 3370.8936
$ws.Range("I136").Value = 2818.7693
$ws.Range("J136").Value = 6062.5
$ws.Range("K136").Value = 8456.3079
$ws.Range("L136").Value = 18187.5
$ws.Range("M136").Value = -5906.3079
$ws.Range("N136").Value = -23287.5
